# Update the requirement/test data: change the sample email addresses
# used for "alireza" (row 2) and "sara" (row 3), and make sure both
# e-mail cells carry a working mailto: hyperlink (row 3 already had one,
# row 2 gets a new one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on B3 so we can recreate it against the
# new address further down (this keeps the relationship id slot so it
# gets reused instead of leaving a stale/duplicate hyperlink behind).
$ws.Range("B3").Hyperlinks.Delete()

# New e-mail values.
$ws.Range("B2").Value = "alireza.13rafe@gmail.com"
$ws.Range("B3").Value = "alireza.00gaming@gmail.com"

# Recreate the hyperlinks for the updated addresses. B3 is added first so
# it reclaims the first relationship id, then B2 gets a new one.
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:alireza.00gaming@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:alireza.13rafe@gmail.com")

# Hyperlinks.Add() assigns a fresh cell style; reapply the built-in
# "Hyperlink" style explicitly so both cells use the same shared style
# (matching the style already used by B3 originally).
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
